$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("F1").Value = "CheckinDate"
$ws.Range("G1").Value = "CheckoutDate"
$ws.Range("F1:G1").Interior.Color = 65535
$ws.Range("F1:G1").Font.Bold = $true

# New column widths
$ws.Columns.Item(6).ColumnWidth = 13.85547
$ws.Columns.Item(7).ColumnWidth = 14.28516

# New date values for row 2
$ws.Range("F2").Value = 45628
$ws.Range("G2").Value = 45537
$ws.Range("F2:G2").NumberFormat = "mm-dd-yy"

# Update selection to match the committed workbook state
$ws.Range("J9").Select()
